$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------------
$headers = @("Name", "email id", "username", "age", "designation", "bio", "work ex ", "Role", "password")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# ---------------------------------------------------------------------------
# Data rows (row index in sheet, username/name, work-ex value)
# ---------------------------------------------------------------------------
$rows = @(
    @{ Row = 2; User = "tb_0"; WorkEx = 5 },
    @{ Row = 3; User = "tb_1"; WorkEx = 6 },
    @{ Row = 4; User = "tb_2"; WorkEx = 7 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $user = $r.User
    $ws.Cells.Item($row, 1).Value = $user          # A: Name
    $ws.Cells.Item($row, 2).Value = "tb@g.com"      # B: email id
    $ws.Cells.Item($row, 3).Value = $user           # C: username
    $ws.Cells.Item($row, 4).Value = 23              # D: age
    $ws.Cells.Item($row, 5).Value = "sdet"          # E: designation
    $ws.Cells.Item($row, 6).Value = "coder"         # F: bio
    $ws.Cells.Item($row, 7).Value = $r.WorkEx       # G: work ex
    $ws.Cells.Item($row, 8).Value = "Employee"      # H: Role
    $ws.Cells.Item($row, 9).Value = "password"      # I: password
}

# ---------------------------------------------------------------------------
# Hyperlinks on the email column
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:tb@g.com")
$ws.Hyperlinks.Add($ws.Range("B3:B4"), "mailto:tb@g.com", "", "", "tb@g.com")
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("B4").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Column widths (approximate the authored widths; engine rounds to the
# nearest 1/6 of a character)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 12.72135
$ws.Columns.Item(2).ColumnWidth = 12.94401
$ws.Columns.Item(3).ColumnWidth = 17.94401
$ws.Columns.Item(4).ColumnWidth = 16.60807
$ws.Columns.Item(5).ColumnWidth = 17.05339
$ws.Columns.Item(6).ColumnWidth = 17.16667
$ws.Columns.Item(7).ColumnWidth = 17.05339
$ws.Columns.Item(8).ColumnWidth = 16.83073
$ws.Columns.Item(9).ColumnWidth = 16.83073

# ---------------------------------------------------------------------------
# Selection, matching the saved view state
# ---------------------------------------------------------------------------
$ws.Range("F8").Select()
